$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three test-case identifiers (corrección de error behavior)
$ws.Range("F2").Value = "'0420172010219"
$ws.Range("F3").Value = "'0420172010219 "
$ws.Range("F4").Value = "'0420194406901 "

# Update current selection to reflect the new active cell
$ws.Range("F7").Select()
